$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'56.384.52"
$ws.Range("E2").Value = '  -1.33%  '

$ws.Range("D3").Value = "'2.970.64"
$ws.Range("E3").Value = '  -3.00%  '

$ws.Range("E4").Value = '  -0.21%  '

$ws.Range("D5").Value = "'501.37"
$ws.Range("E5").Value = '  -2.01%  '

$ws.Range("D6").Value = "'134.17"
$ws.Range("E6").Value = '  +5.26%  '

$ws.Range("E7").Value = '  -0.19%  '

$ws.Range("D8").Value = "'0.428"
$ws.Range("E8").Value = '  -0.59%  '

$ws.Range("D9").Value = "'7.26"
$ws.Range("E9").Value = '  +2.77%  '

$ws.Range("D10").Value = "'0.106"
$ws.Range("E10").Value = '  +2.81%  '

$ws.Range("D11").Value = "'0.350"
$ws.Range("E11").Value = '  -2.46%  '

$ws.Range("D13").Value = "'3.481.41"
$ws.Range("E13").Value = '  -3.86%  '

$ws.Range("D14").Value = "'25.09"
$ws.Range("E14").Value = '  +3.77%  '

$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = "'56.377.27"
$ws.Range("E15").Value = '  +3.12%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = "'0.0000150"
$ws.Range("E16").Value = '  +3.99%  '

$ws.Range("D17").Value = "'2.974.77"
$ws.Range("E17").Value = '  -3.72%  '

$ws.Range("D18").Value = "'5.67"
$ws.Range("E18").Value = '  +3.06%  '

$ws.Range("D19").Value = "'12.29"
$ws.Range("E19").Value = '  -0.96%  '

$ws.Range("D20").Value = "'7.76"
$ws.Range("E20").Value = '  +2.78%  '

$ws.Range("D21").Value = "'327.45"
$ws.Range("E21").Value = '  +0.00%  '

$ws.Range("E22").Value = '  +0.16%  '

$ws.Range("D23").Value = "'0.468"
$ws.Range("E23").Value = '  -4.09%  '

$ws.Range("D24").Value = "'62.12"
$ws.Range("E24").Value = '  -5.10%  '

$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = '  +0.06%  '

$ws.Range("D26").Value = "'0.163"
$ws.Range("E26").Value = '  -1.81%  '

$ws.Range("D27").Value = "'0.0₃0895"
$ws.Range("E27").Value = '  +2.67%  '

$ws.Range("D28").Value = "'0.998"
$ws.Range("E28").Value = '  -0.12%  '

$ws.Range("D29").Value = "'6.45"
$ws.Range("E29").Value = '  -0.78%  '

$ws.Range("D30").Value = "'6.77"
$ws.Range("E30").Value = '  +2.94%  '

$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = "'1.74"
$ws.Range("E31").Value = '  -2.52%  '

$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").Value = "'1.16"
$ws.Range("E32").Value = '  -4.38%  '

$ws.Range("D33").Value = "'20.24"
$ws.Range("E33").Value = '  -1.78%  '

$ws.Range("D34").Value = "'155.14"
$ws.Range("E34").Value = '  -0.64%  '

$ws.Range("D35").Value = "'4.43"
$ws.Range("E35").Value = '  -2.77%  '

$ws.Range("D36").Value = "'1.28"
$ws.Range("E36").Value = '  -1.63%  '

$ws.Range("D37").Value = "'5.53"
$ws.Range("E37").Value = '  -5.79%  '

$ws.Range("D38").Value = "'0.0670"
$ws.Range("E38").Value = '  +2.37%  '

$ws.Range("D39").Value = "'22.86"
$ws.Range("E39").Value = '  +0.00%  '

$ws.Range("D40").Value = "'3.005.47"
$ws.Range("E40").Value = '  -3.43%  '

$ws.Range("E41").Value = '  -0.09%  '

$ws.Range("D42").Value = "'35.81"
$ws.Range("E42").Value = '  -1.04%  '

$ws.Range("D43").Value = "'0.639"
$ws.Range("E43").Value = '  -3.58%  '

$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = "'2.230.44"
$ws.Range("E44").Value = '  +0.52%  '

$ws.Range("B45").Value = 'ONDO'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D45").Value = "'0.989"
$ws.Range("E45").Value = '  -4.44%  '

$ws.Range("D46").Value = "'1.38"
$ws.Range("E46").Value = '  +1.21%  '

$ws.Range("D47").Value = "'3.53"
$ws.Range("E47").Value = '  -4.59%  '

$ws.Range("D48").Value = "'1.94"
$ws.Range("E48").Value = '  +13.72%  '

$ws.Range("D49").Value = "'0.0234"
$ws.Range("E49").Value = '  +3.78%  '

$ws.Range("D50").Value = "'5.75"
$ws.Range("E50").Value = '  -3.03%  '

$ws.Range("D51").Value = "'18.88"
$ws.Range("E51").Value = '  -3.21%  '
